# Applies the "HTML Output files checkin" edit:
#  - ProcessPayrollForWeeklyTax (sheet3) rows 2-6: the per-employer "which
#    employer" text that used to live in column C moves left into column B
#    (replacing the constant "DONT TOUCH AUTO 2040W EMPLOYER" placeholder);
#    the "2040W_Payroll" text that used to live in column D moves left into
#    column C (losing the old highlighted/hyperlink-ish style); and a new
#    "Weekly" value is written into column D, formatted like column E.
#  - The stale hyperlink that used to sit over B2:B6 is removed (it no
#    longer makes sense once B holds plain employer-name text).
#  - Leftover view state (scrolled-off topLeftCell / odd selections) on the
#    ProcessPayrollForWeeklyTax and TestReports sheets is reset.

$wb = $excel.ActiveWorkbook

# ---- ProcessPayrollForWeeklyTax --------------------------------------
$ws3 = $wb.Worksheets.Item("ProcessPayrollForWeeklyTax")

for ($r = 2; $r -le 6; $r++) {
    $colB = $ws3.Range("B$r")
    $colC = $ws3.Range("C$r")
    $colD = $ws3.Range("D$r")
    $colE = $ws3.Range("E$r")

    # Shift the real content two columns left (B keeps its existing style,
    # which already matches the incoming value's old style).
    $colB.Value2 = $colC.Value2
    $colC.Value2 = $colD.Value2
    $colC.Style = "Normal"

    # New column D value, formatted the same way as column E.
    $colD.Value2 = "Weekly"
    $colE.Copy()
    $colD.PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# The hyperlink over the old column B no longer applies.
while ($ws3.Hyperlinks.Count -gt 0) {
    $ws3.Hyperlinks.Item(1).Delete()
}

# Reset the scrolled/odd selection left over from editing.
$ws3.Activate()
$ws3.Range("C9").Select()

# ---- TestReports -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("TestReports")
$ws4.Activate()
$ws4.Range("D2").Select()

# Leave the originally active sheet selected.
$ws3.Activate()
